$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 5.25
$ws.Range("L2").Value = 5.5
$ws.Range("Q2").Value = 1.93
$ws.Range("R2").Value = 1.93
$ws.Range("V2").Value = 1.24
$ws.Range("AA2").Value = 2.25
$ws.Range("AB2").Value = 1.57
$ws.Range("AG2").Value = 19
$ws.Range("AN2").Value = 10
$ws.Range("AO2").Value = 23
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.88
$ws.Range("S3").Value = 2.6
$ws.Range("T3").Value = 1.48
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("I7").Value = 3.25
$ws.Range("K7").Value = 1.83
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.78
$ws.Range("G9").Value = 2.35
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("AA9").Value = 1.91
$ws.Range("AB9").Value = 1.91
$ws.Range("AE9").Value = 9.5
$ws.Range("G10").Value = 3.7
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 1.95
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 10
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 1.7
$ws.Range("W10").Value = 3.75
$ws.Range("X10").Value = 1.25
$ws.Range("Y10").Value = 1.44
$ws.Range("Z10").Value = 2.63
$ws.Range("AD10").Value = 19
$ws.Range("AO10").Value = 8.5
$ws.Range("AS10").Value = 29
$ws.Range("J13").Value = 1.91
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 8
$ws.Range("AA14").Value = 1.8
$ws.Range("AB14").Value = 1.91
